$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns F, G, H
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Match the header style used by the existing header row (A1:E1)
$ws.Range("F1").Style = $ws.Range("E1").Style
$ws.Range("G1").Style = $ws.Range("E1").Style
$ws.Range("H1").Style = $ws.Range("E1").Style

# Boolean data for rows 2-8
$values = @(
    @($false, $true,  $false),
    @($false, $false, $false),
    @($false, $false, $false),
    @($false, $true,  $false),
    @($false, $true,  $false),
    @($false, $false, $false),
    @($false, $false, $false)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Range("F$row").Value = $values[$i][0]
    $ws.Range("G$row").Value = $values[$i][1]
    $ws.Range("H$row").Value = $values[$i][2]
}
